# Apply the "Add descriptions titles" edit to the DMI Identifiant Local
# Distributeur StructureDefinition workbook.
#
# Changes:
#   Metadata sheet:
#     - Title       (B5)  : set to "DMI Identifiant Local Distributeur"
#     - Date        (B8)  : updated to "2026-02-25T08:15:31+00:00"
#     - Description (B12) : set to the French description text
#   Elements sheet (row 2 = root "Extension" element):
#     - Short (L2)               : "DMI Identifiant Local Distributeur"
#     - Definition (M2)          : the French description text
#     - Mapping: RIM Mapping (AK2): cleared (was "N/A")

$wb = $excel.ActiveWorkbook

$titleText = "DMI Identifiant Local Distributeur"
$descriptionText = "Extension créée dans ce volet pour représenter l'identifiant local distributeur."
$dateText = "2026-02-25T08:15:31+00:00"

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B5").Value = $titleText
$metadata.Range("B8").Value = $dateText
$metadata.Range("B12").Value = $descriptionText

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("L2").Value = $titleText
$elements.Range("M2").Value = $descriptionText
$elements.Range("AK2").Value = ""
